$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: "Sending cluster" changes from MuSCs to ECs, with refreshed TPM values ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pmch"
$ws.Range("C2").Value = "Mchr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.2423523333333333
$ws.Range("H2").Value = 0.727057
$ws.Range("I2").Value = 0.7762337358684754
$ws.Range("J2").Value = 0.7762337358684754
$ws.Range("M2").Value = 0.3559683333333334
$ws.Range("Q2").Value = 0.08626975617611111
$ws.Range("R2").Value = 0.776427805585
$ws.Range("S2").Value = 0.7762337358684754
$ws.Range("T2").Value = 0.7762337358684754

# --- Add new row 3 for "MuSCs" sending cluster with its own TPM values ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Pmch"
$ws.Range("C3").Value = "Mchr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06986333333333333
$ws.Range("H3").Value = 0.20959
$ws.Range("I3").Value = 0.2237662641315245
$ws.Range("J3").Value = 0.2237662641315245
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.3559683333333334
$ws.Range("N3").Value = 1.067905
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.02486913432777778
$ws.Range("R3").Value = 0.22382220895
$ws.Range("S3").Value = 0.2237662641315245
$ws.Range("T3").Value = 0.2237662641315245
